$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 2.2.0-ballot -> 2.1.0
$meta.Range("B3").Value = "2.1.0"

# Date: 2025-12-19T08:32:44+00:00 -> 2025-12-19T08:44:55+00:00
$meta.Range("B8").Value = "2025-12-19T08:44:55+00:00"

# Base Definition: drop the "|4.0.1" version suffix
$meta.Range("B18").Value = "http://hl7.org/fhir/StructureDefinition/Extension"

# --- Elements sheet updates ---
$elements = $wb.Worksheets.Item("Elements")

# Extension.value[x] Type(s): drop the "|2.2.0-ballot" version suffix from the Reference(...) target,
# keeping the trailing newline that was already part of the cell text.
$elements.Range("K6").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/tddui/StructureDefinition/tddui-document-reference)`n"

# Column K width: 89.21875 -> 79.70703125 (character units).
# The host snaps ColumnWidth to whole-pixel boundaries, so the input value
# is chosen such that the stored/rounded result lands as close as possible
# to the target width (79.666... is the nearest representable value).
$elements.Columns.Item(11).ColumnWidth = 78.8
